$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 58 fix: AD58/AE58/AF58 were "4" and should be "0"
# ---------------------------------------------------------------------------
$ws.Range("AD58").Value = "0"
$ws.Range("AE58").Value = "0"
$ws.Range("AF58").Value = "0"

# ---------------------------------------------------------------------------
# Columns A:AC (and AA:AC) share the same value across the four new rows
# (59-62), mirroring row 58. Only AD:AF (viterations/piterations/diterations)
# differ per-row, running through the sequence 1,2,3,4.
# ---------------------------------------------------------------------------
$commonValues = @{
    "A"  = "/Users/davl3232/Documents/uni/tg/Ane-stent/modelos/cilindroRadio.vtk"
    "B"  = "0"
    "C"  = "1"
    "D"  = "1"
    "E"  = "1"
    "F"  = "90"
    "G"  = "0"
    "H"  = "0"
    "I"  = "0"
    "J"  = "4"
    "K"  = "0"
    "L"  = "0.1"
    "M"  = "0.1"
    "N"  = "0.1"
    "O"  = "10"
    "P"  = "1"
    "V"  = "1"
    "AA" = "0.1"
    "AB" = "1"
    "AC" = "1"
}

# Columns that must be stored as real numbers (no text coercion), matching
# the original numeric-typed cells already present in row 58.
$numericColumns = @("Q", "R", "S", "T", "U", "W", "X", "Y", "Z")
$numericValues = @{
    "Q" = 0
    "R" = 0
    "S" = 0
    "T" = 0
    "U" = 0
    "W" = 0
    "X" = 1
    "Y" = 1
    "Z" = 1
}

$newRows = @(59, 60, 61, 62)
$adValues = @{ 59 = "1"; 60 = "2"; 61 = "3"; 62 = "4" }

foreach ($row in $newRows) {
    foreach ($col in $commonValues.Keys) {
        $ws.Range("$col$row").Value = $commonValues[$col]
    }

    foreach ($col in $numericColumns) {
        $cell = $ws.Range("$col$row")
        # Every column in this sheet is pre-formatted as Text ("@"), which
        # forces any assigned value to be stored as a string. Temporarily
        # switching to General preserves the numeric storage of the value,
        # then restoring "@" keeps the visual formatting identical to the
        # rest of the sheet.
        $cell.NumberFormat = "General"
        $cell.Value = $numericValues[$col]
        $cell.NumberFormat = "@"
    }

    $ad = $adValues[$row]
    $ws.Range("AD$row").Value = $ad
    $ws.Range("AE$row").Value = $ad
    $ws.Range("AF$row").Value = $ad
}

# ---------------------------------------------------------------------------
# View state: scroll back to the left (drop the U1 top-left override) and
# move the active selection down to the freshly added last row.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("A62:AF62").Select()
$excel.ActiveWindow.RangeSelection.Item(1).Activate()
